# Apply updated "dSF" (column F) values as part of a repull/recalculation
# of the farmer_buck data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -10
    5  = -1
    6  = -1
    7  = -7
    15 = -1
    21 = -2
    22 = 0
    28 = 0
    30 = -3
    34 = -4
    35 = 0
    37 = -2
    40 = 0
    51 = 2
    54 = -2
    55 = 1
    59 = -3
    60 = 0
    65 = 4
    67 = -1
    72 = -1
    78 = -3
    79 = 1
    80 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
